$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update student first name (Imię ucznia) and last name (Nazwisko ucznia)
$ws.Range("C2").Value = "Olek"
$ws.Range("D2").Value = "Ura"

# Move selection as a natural side-effect of editing, matching the diff
$ws.Range("C6").Select()
